$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.02354566666666667
$ws.Range("H2").Value = 0.07063700000000001
$ws.Range("I2").Value = 0.002815555392485919
$ws.Range("J2").Value = 0.002815555392485918
$ws.Range("M2").Value = 23.28617366666667
$ws.Range("N2").Value = 69.858521
$ws.Range("O2").Value = 0.2304887056246027
$ws.Range("P2").Value = 0.2304887056246027
$ws.Range("Q2").Value = 0.5482884830974445
$ws.Range("R2").Value = 4.934596347877
$ws.Range("S2").Value = 0.0006489537180284497
$ws.Range("T2").Value = 0.0006489537180284495
$ws.Range("G3").Value = 0.02354566666666667
$ws.Range("H3").Value = 0.07063700000000001
$ws.Range("I3").Value = 0.002815555392485919
$ws.Range("J3").Value = 0.002815555392485918
$ws.Range("O3").Value = 0.007098179626924059
$ws.Range("P3").Value = 0.007098179626924059
$ws.Range("Q3").Value = 0.016885209754
$ws.Range("R3").Value = 0.151966887786
$ws.Range("S3").Value = 0.00001998531792541972
$ws.Range("T3").Value = 0.00001998531792541972
$ws.Range("G4").Value = 0.02354566666666667
$ws.Range("H4").Value = 0.07063700000000001
$ws.Range("I4").Value = 0.002815555392485919
$ws.Range("J4").Value = 0.002815555392485918
$ws.Range("O4").Value = 0.7624131147484733
$ws.Range("P4").Value = 0.7624131147484732
$ws.Range("Q4").Value = 1.813634768116889
$ws.Range("R4").Value = 16.322712913052
$ws.Range("S4").Value = 0.00214661635653205
$ws.Range("T4").Value = 0.002146616356532049
$ws.Range("I5").Value = 0.9868456480383168
$ws.Range("J5").Value = 0.9868456480383166
$ws.Range("M5").Value = 23.28617366666667
$ws.Range("N5").Value = 69.858521
$ws.Range("O5").Value = 0.2304887056246027
$ws.Range("P5").Value = 0.2304887056246027
$ws.Range("Q5").Value = 192.1738442291895
$ws.Range("R5").Value = 1729.564598062705
$ws.Range("S5").Value = 0.2274567760676239
$ws.Range("T5").Value = 0.2274567760676239
$ws.Range("I6").Value = 0.9868456480383168
$ws.Range("J6").Value = 0.9868456480383166
$ws.Range("O6").Value = 0.007098179626924059
$ws.Range("P6").Value = 0.007098179626924059
$ws.Range("S6").Value = 0.007004807673824251
$ws.Range("T6").Value = 0.00700480767382425
$ws.Range("I7").Value = 0.9868456480383168
$ws.Range("J7").Value = 0.9868456480383166
$ws.Range("O7").Value = 0.7624131147484733
$ws.Range("P7").Value = 0.7624131147484732
$ws.Range("S7").Value = 0.7523840642968687
$ws.Range("T7").Value = 0.7523840642968684
$ws.Range("I8").Value = 0.0103387965691973
$ws.Range("J8").Value = 0.0103387965691973
$ws.Range("M8").Value = 23.28617366666667
$ws.Range("N8").Value = 69.858521
$ws.Range("O8").Value = 0.2304887056246027
$ws.Range("P8").Value = 0.2304887056246027
$ws.Range("Q8").Value = 2.013330337277889
$ws.Range("R8").Value = 18.119973035501
$ws.Range("S8").Value = 0.00238297583895037
$ws.Range("T8").Value = 0.00238297583895037
$ws.Range("I9").Value = 0.0103387965691973
$ws.Range("J9").Value = 0.0103387965691973
$ws.Range("O9").Value = 0.007098179626924059
$ws.Range("P9").Value = 0.007098179626924059
$ws.Range("S9").Value = 0.00007338663517438865
$ws.Range("T9").Value = 0.00007338663517438865
$ws.Range("I10").Value = 0.0103387965691973
$ws.Range("J10").Value = 0.0103387965691973
$ws.Range("O10").Value = 0.7624131147484733
$ws.Range("P10").Value = 0.7624131147484732
$ws.Range("S10").Value = 0.007882434095072546
$ws.Range("T10").Value = 0.007882434095072544
